$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking price strings stay as text (matching the
# original inline-string cells) instead of being auto-converted to numbers.
# (NumberFormat must be set per-cell; multi-area ranges only apply to the
# first area in this runtime.)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.759.89"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "2.498.71"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "587.82"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").Value = "176.28"
$ws.Range("E6").Value = "  +4.04%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "0.517"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("E9").Value = "  +5.98%  "
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("E11").Value = "  +4.01%  "
$ws.Range("D12").Value = "4.95"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").Value = "25.85"
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").Value = "2.917.25"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "67.547.82"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").Value = "0.0000173"
$ws.Range("E16").Value = "  +2.73%  "
$ws.Range("D17").Value = "2.514.20"
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("D18").Value = "11.12"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("D19").Value = "7.56"
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("D20").Value = "352.24"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").Value = "4.10"
$ws.Range("E21").Value = "  +2.42%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "70.56"
$ws.Range("E23").Value = "  +3.18%  "
$ws.Range("D24").Value = "4.31"
$ws.Range("E24").Value = "  +3.05%  "
$ws.Range("D25").Value = "1.77"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("D26").Value = "9.20"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("D27").Value = "2.624.52"
$ws.Range("E27").Value = "  +1.50%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "0.0₃0916"
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("D30").Value = "510.14"
$ws.Range("E30").Value = "  +0.74%  "
$ws.Range("D31").Value = "7.85"
$ws.Range("E31").Value = "  +3.39%  "
$ws.Range("E32").Value = "  +3.86%  "
$ws.Range("E33").Value = "  +1.52%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  +8.45%  "
$ws.Range("D36").Value = "163.78"
$ws.Range("E36").Value = "  +3.31%  "
$ws.Range("D37").Value = "18.51"
$ws.Range("E37").Value = "  +1.87%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +4.77%  "
$ws.Range("E42").Value = "  +1.98%  "
$ws.Range("D43").Value = "4.90"
$ws.Range("E43").Value = "  +2.27%  "
$ws.Range("D44").Value = "2.43"
$ws.Range("E44").Value = "  +3.66%  "
$ws.Range("D45").Value = "146.31"
$ws.Range("E45").Value = "  +3.86%  "
$ws.Range("D46").Value = "3.53"
$ws.Range("E46").Value = "  +3.15%  "
$ws.Range("D47").Value = "0.518"
$ws.Range("E47").Value = "  +1.79%  "
$ws.Range("D48").Value = "0.0₆0258"
$ws.Range("E48").Value = "  +3.39%  "
$ws.Range("D49").Value = "0.0747"
$ws.Range("E49").Value = "  +2.70%  "
$ws.Range("D50").Value = "1.60"
$ws.Range("E50").Value = "  +2.39%  "
$ws.Range("D51").Value = "0.588"
$ws.Range("E51").Value = "  +1.40%  "

